# Add three new R packages (ggpacman, dgapaq, dmapaq) to the "packages" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("packages")
$ws.Activate()

# Row 11 (ggpacman) is entered first.
$ws.Range("A11").Value = "ggpacman"
$ws.Range("B11").Value = "A 'ggplot2' and 'gganimate' Version of Pac-Man"
$ws.Range("C11").Value = "Mickaël Canouil"
$ws.Range("D11").Value = "Author / Creator"
$ws.Range("E11").Value = "BOTH"
$ws.Range("F11").Value = "mcanouil"
$ws.Range("G11").Value = "2020-05"

# Then the "since" column is filled down for the two new rows 9 and 10.
$ws.Range("G9").Value = "2020-03"
$ws.Range("G10").Value = "2020-03"

# Row 9 (dgapaq) name/title.
$ws.Range("B9").Value = "DNA Genotyping Arrays Processing And Quality-Control"
$ws.Range("A9").Value = "dgapaq"

# Row 10 (dmapaq) name/title.
$ws.Range("A10").Value = "dmapaq"
$ws.Range("B10").Value = "DNA Methylation Arrays Processing And Quality-Control"

# Remaining cells for rows 9 and 10.
$ws.Range("C9").Value = "Mickaël Canouil and Lijiao Ning"
$ws.Range("D9").Value = "Author / Creator"
$ws.Range("E9").Value = "GitHub"
$ws.Range("F9").Value = "omicsr"

$ws.Range("C10").Value = "Mickaël Canouil, Mathilde Boissel and Lijiao Ning"
$ws.Range("D10").Value = "Author / Creator"
$ws.Range("E10").Value = "GitHub"
$ws.Range("F10").Value = "omicsr"

# Widen column B (title) to fit the new, longer entries.
$ws.Columns("B:B").ColumnWidth = 50.45

# Page setup for the sheet (paper size / orientation).
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Final UI selection left on the sheet.
$ws.Range("D21").Select()
